$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @(44319, 3, 31, 180.9162532827546),
    @(44320, 1, 30, 175.0802451123432),
    @(44321, 3, 32, 186.752261453166)
)

$startRow = 245
$lastRow = $startRow - 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Copy formatting from the cell directly above so the new row matches
    # the existing style (date style with border/centering).
    $ws.Range("A" + $lastRow).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$excel.CutCopyMode = 0
